$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value2 = 0.04997382157907282
$ws.Range("D2").Value2 = 0.01386200908951452
$ws.Range("E2").Value2 = 0.4243101376564908
$ws.Range("F2").Value2 = 0.4515580344483539
$ws.Range("G2").Value2 = 0.3013419787874909
$ws.Range("H2").Value2 = 0.4455450221618378
$ws.Range("K2").Value2 = 1.077123396043532
$ws.Range("N2").Value2 = 0.8880294875536237
$ws.Range("O2").Value2 = 1.420488363125202

$ws.Range("C3").Value2 = 0.04434312488731962
$ws.Range("D3").Value2 = 0.01228864329191026
$ws.Range("E3").Value2 = 0.3701214267360058
$ws.Range("F3").Value2 = 0.4443250307711821
$ws.Range("G3").Value2 = 0.2951115846810808
$ws.Range("H3").Value2 = 0.4468836137525045
$ws.Range("K3").Value2 = 0.9400804986817093
$ws.Range("N3").Value2 = 0.8839194764473888
$ws.Range("O3").Value2 = 1.409775580676836

$ws.Range("C4").Value2 = 0.04089948348679684
$ws.Range("D4").Value2 = 0.01131724892957919
$ws.Range("E4").Value2 = 0.3369378955857343
$ws.Range("F4").Value2 = 0.4402722654424593
$ws.Range("G4").Value2 = 0.2916004550931746
$ws.Range("H4").Value2 = 0.4479939264975314
$ws.Range("K4").Value2 = 0.8556671847902919
$ws.Range("N4").Value2 = 0.8817473418291257
$ws.Range("O4").Value2 = 1.404457330437936

$ws.Range("C5").Value2 = 0.03949959210599729
$ws.Range("D5").Value2 = 0.01092008695893298
$ws.Range("E5").Value2 = 0.3234358543398628
$ws.Range("F5").Value2 = 0.438718046363789
$ws.Range("G5").Value2 = 0.290248233801961
$ws.Range("H5").Value2 = 0.4485187927989358
$ws.Range("K5").Value2 = 0.8212024614053917
$ws.Range("N5").Value2 = 0.880950876720334
$ws.Range("O5").Value2 = 1.402605575881125

$ws.Range("C6").Value2 = 0.03926734760882766
$ws.Range("D6").Value2 = 0.0108540602945979
$ws.Range("E6").Value2 = 0.3211950388773772
$ws.Range("F6").Value2 = 0.4384658373836885
$ws.Range("G6").Value2 = 0.2900284319685937
$ws.Range("H6").Value2 = 0.4486103160648867
$ws.Range("K6").Value2 = 0.8154757131872259
$ws.Range("N6").Value2 = 0.8808239911864462
$ws.Range("O6").Value2 = 1.402317111713842

$ws.Range("C7").Value2 = 0.04088059022112134
$ws.Range("D7").Value2 = 0.01131189792906184
$ws.Range("E7").Value2 = 0.336755721559058
$ws.Range("F7").Value2 = 0.4402509111153208
$ws.Range("G7").Value2 = 0.2915819009499643
$ws.Range("H7").Value2 = 0.4480007120232301
$ws.Range("K7").Value2 = 0.8552026446357104
$ws.Range("N7").Value2 = 0.8817362408525611
$ws.Range("O7").Value2 = 1.404431081305233

$ws.Range("C8").Value2 = 0.04802952286730999
$ws.Range("D8").Value2 = 0.01332064000211375
$ws.Range("E8").Value2 = 0.4056064188250019
$ws.Range("F8").Value2 = 0.4489832933986406
$ws.Range("G8").Value2 = 0.29912821717447
$ws.Range("H8").Value2 = 0.4459466360289071
$ws.Range("K8").Value2 = 1.02992780217744
$ws.Range("N8").Value2 = 0.8865396076148784
$ws.Range("O8").Value2 = 1.416532455323789

$ws.Range("C9").Value2 = 0.06215772947155074
$ws.Range("D9").Value2 = 0.01721616289282224
$ws.Range("E9").Value2 = 0.5414103827941688
$ws.Range("F9").Value2 = 0.4692053025033971
$ws.Range("G9").Value2 = 0.3164425627586525
$ws.Range("H9").Value2 = 0.4442128710571325
$ws.Range("K9").Value2 = 1.370367210366453
$ws.Range("N9").Value2 = 0.8987347429896744
$ws.Range("O9").Value2 = 1.450315031536377

$ws.Range("C10").Value2 = 0.07260675156352647
$ws.Range("D10").Value2 = 0.02005023012586094
$ws.Range("E10").Value2 = 0.6417931197073585
$ws.Range("F10").Value2 = 0.4859762851314997
$ws.Range("G10").Value2 = 0.3307291313723653
$ws.Range("H10").Value2 = 0.4443465404101374
$ws.Range("K10").Value2 = 1.619084449520187
$ws.Range("N10").Value2 = 0.9093721544169853
$ws.Range("O10").Value2 = 1.481349453971802

$ws.Range("C11").Value2 = 0.07737586259557361
$ws.Range("D11").Value2 = 0.02133317191321282
$ws.Range("E11").Value2 = 0.6876198744579511
$ws.Range("F11").Value2 = 0.4940270216706608
$ws.Range("G11").Value2 = 0.3375754111214775
$ws.Range("H11").Value2 = 0.4447149295095869
$ws.Range("K11").Value2 = 1.731915469411433
$ws.Range("N11").Value2 = 0.9145728895110778
$ws.Range("O11").Value2 = 1.496836188275864

$ws.Range("C12").Value2 = 0.07918409612301502
$ws.Range("D12").Value2 = 0.02181805543579429
$ws.Range("E12").Value2 = 0.7049986939847486
$ws.Range("F12").Value2 = 0.4971366475163421
$ws.Range("G12").Value2 = 0.3402183806311996
$ws.Range("H12").Value2 = 0.4448988082245506
$ws.Range("K12").Value2 = 1.774595267241921
$ws.Range("N12").Value2 = 0.9165940287818017
$ws.Range("O12").Value2 = 1.502898914123051

$ws.Range("C13").Value2 = 0.07879455946793712
$ws.Range("D13").Value2 = 0.02171366938468111
$ws.Range("E13").Value2 = 0.7012547019538431
$ws.Range("F13").Value2 = 0.4964642156444654
$ws.Range("G13").Value2 = 0.3396469189438562
$ws.Range("H13").Value2 = 0.4448572305759626
$ws.Range("K13").Value2 = 1.765405525608628
$ws.Range("N13").Value2 = 0.916156443799764
$ws.Range("O13").Value2 = 1.501584359159892

$ws.Range("C14").Value2 = 0.07752458145250785
$ws.Range("D14").Value2 = 0.02137308258095061
$ws.Range("E14").Value2 = 0.6890491239283278
$ws.Range("F14").Value2 = 0.4942816280717039
$ws.Range("G14").Value2 = 0.3377918362267849
$ws.Range("H14").Value2 = 0.4447291670758915
$ws.Range("K14").Value2 = 1.735427712856563
$ws.Range("N14").Value2 = 0.914738134617437
$ws.Range("O14").Value2 = 1.497330991569726

$ws.Range("C15").Value2 = 0.07674697956029775
$ws.Range("D15").Value2 = 0.02116434007111678
$ws.Range("E15").Value2 = 0.6815761964983267
$ws.Range("F15").Value2 = 0.4929526834365276
$ws.Range("G15").Value2 = 0.3366621270220094
$ws.Range("H15").Value2 = 0.4446565081809553
$ws.Range("K15").Value2 = 1.717059274018538
$ws.Range("N15").Value2 = 0.9138761083359839
$ws.Range("O15").Value2 = 1.494751538028567

$ws.Range("C16").Value2 = 0.07229539709884136
$ws.Range("D16").Value2 = 0.01996625736682489
$ws.Range("E16").Value2 = 0.638801658778732
$ws.Range("F16").Value2 = 0.4854586645604257
$ws.Range("G16").Value2 = 0.3302887409712696
$ws.Range("H16").Value2 = 0.4443286672704545
$ws.Range("K16").Value2 = 1.611704230276416
$ws.Range("N16").Value2 = 0.9090395286832802
$ws.Range("O16").Value2 = 1.480365024872128

$ws.Range("C17").Value2 = 0.06956854391917489
$ws.Range("D17").Value2 = 0.01922963707136205
$ws.Range("E17").Value2 = 0.6126037447600794
$ws.Range("F17").Value2 = 0.4809695549568858
$ws.Range("G17").Value2 = 0.3264681678130756
$ws.Range("H17").Value2 = 0.4442064315573759
$ws.Range("K17").Value2 = 1.546991035250073
$ws.Range("N17").Value2 = 0.9061648867451595
$ws.Range("O17").Value2 = 1.471890979218642

$ws.Range("C18").Value2 = 0.06800161696108376
$ws.Range("D18").Value2 = 0.01880536258514098
$ws.Range("E18").Value2 = 0.5975505443285982
$ws.Range("F18").Value2 = 0.4784271922708427
$ws.Range("G18").Value2 = 0.3243033299875577
$ws.Range("H18").Value2 = 0.4441650670010375
$ws.Range("K18").Value2 = 1.509740499519012
$ws.Range("N18").Value2 = 0.9045455321361402
$ws.Range("O18").Value2 = 1.467145708648019

$ws.Range("C19").Value2 = 0.06747133769600566
$ws.Range("D19").Value2 = 0.01866161025331792
$ws.Range("E19").Value2 = 0.5924563329920858
$ws.Range("F19").Value2 = 0.477573191347588
$ws.Range("G19").Value2 = 0.3235759464272405
$ws.Range("H19").Value2 = 0.4441560280373267
$ws.Range("K19").Value2 = 1.497123152345807
$ws.Range("N19").Value2 = 0.904003105632313
$ws.Range("O19").Value2 = 1.465561117284693

$ws.Range("C20").Value2 = 0.06985866818980924
$ws.Range("D20").Value2 = 0.01930811286872114
$ws.Range("E20").Value2 = 0.6153909760641625
$ws.Range("F20").Value2 = 0.4814433215672125
$ws.Range("G20").Value2 = 0.3268714910912678
$ws.Range("H20").Value2 = 0.4442164472790466
$ws.Range("K20").Value2 = 1.553882904049431
$ws.Range("N20").Value2 = 0.9064673739154898
$ws.Range("O20").Value2 = 1.472779718037003

$ws.Range("C21").Value2 = 0.07789754302339702
$ws.Range("D21").Value2 = 0.02147314692875568
$ws.Range("E21").Value2 = 0.6926334982629072
$ws.Range("F21").Value2 = 0.4949210488309461
$ws.Range("G21").Value2 = 0.3383353463066214
$ws.Range("H21").Value2 = 0.4447655768477006
$ws.Range("K21").Value2 = 1.744234211106061
$ws.Range("N21").Value2 = 0.9151533245684647
$ws.Range("O21").Value2 = 1.498574917720958

$ws.Range("C22").Value2 = 0.08316467225213842
$ws.Range("D22").Value2 = 0.02288263470950369
$ws.Range("E22").Value2 = 0.7432638629730803
$ws.Range("F22").Value2 = 0.50408517198872
$ws.Range("G22").Value2 = 0.3461218260892025
$ws.Range("H22").Value2 = 0.4453832015709906
$ws.Range("K22").Value2 = 1.868365716984101
$ws.Range("N22").Value2 = 0.9211315145176258
$ws.Range("O22").Value2 = 1.516589680422072

$ws.Range("C23").Value2 = 0.08035229207403916
$ws.Range("D23").Value2 = 0.02213087788602763
$ws.Range("E23").Value2 = 0.7162273160135157
$ws.Range("F23").Value2 = 0.4991614424140494
$ws.Range("G23").Value2 = 0.3419389523257621
$ws.Range("H23").Value2 = 0.4450298401189201
$ws.Range("K23").Value2 = 1.802140144629959
$ws.Range("N23").Value2 = 0.9179133530880108
$ws.Range("O23").Value2 = 1.506868624214746

$ws.Range("C24").Value2 = 0.0697275005869642
$ws.Range("D24").Value2 = 0.01927263639297649
$ws.Range("E24").Value2 = 0.6141308430277093
$ws.Range("F24").Value2 = 0.4812290118260165
$ws.Range("G24").Value2 = 0.3266890500523658
$ws.Range("H24").Value2 = 0.4442118291239154
$ws.Range("K24").Value2 = 1.550767233480826
$ws.Range("N24").Value2 = 0.9063305156641803
$ws.Range("O24").Value2 = 1.472377525046653

$ws.Range("C25").Value2 = 0.05832373131028135
$ws.Range("D25").Value2 = 0.01616713397754665
$ws.Range("E25").Value2 = 0.5045743668315197
$ws.Range("F25").Value2 = 0.463400450432907
$ws.Range("G25").Value2 = 0.3114859220056729
$ws.Range("H25").Value2 = 0.4444353750269983
$ws.Range("K25").Value2 = 1.278511180746307
$ws.Range("N25").Value2 = 0.8951400068560247
$ws.Range("O25").Value2 = 1.440090904297506

